$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 37: Date, Tasks done, hours
$ws.Range("A37").Value = 45672
$ws.Range("A37").NumberFormat = $ws.Range("A36").NumberFormat
$ws.Range("B37").Value = "reworking plan for story"
$ws.Range("C37").Value = 5

# Update the view to match the new selection/scroll position
$ws.Range("C38").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
